$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 0.04847177082993426
$ws.Range("D2").Value = 0.1534718802712476
$ws.Range("E2").Value = 0.1364174140543746
$ws.Range("F2").Value = 1.52878314466372
$ws.Range("G2").Value = 0.002472927089828936
$ws.Range("J2").Value = 0.1577116572115003
$ws.Range("K2").Value = 1.376394238103046
$ws.Range("M2").Value = 0.4331660064745222
$ws.Range("O2").Value = 3.812364515691229
$ws.Range("C3").Value = 0.04308372056300414
$ws.Range("D3").Value = 0.1511756676351581
$ws.Range("E3").Value = 0.1364779707666841
$ws.Range("F3").Value = 1.53966367904944
$ws.Range("G3").Value = 0.002476190568244493
$ws.Range("J3").Value = 0.1591681047621876
$ws.Range("K3").Value = 1.22834319122228
$ws.Range("M3").Value = 0.4031164499780644
$ws.Range("O3").Value = 3.850359260727288
$ws.Range("C4").Value = 0.03979032310876107
$ws.Range("D4").Value = 0.1498065855606114
$ws.Range("E4").Value = 0.1365755322531648
$ws.Range("F4").Value = 1.547393368687544
$ws.Range("G4").Value = 0.002478300271214605
$ws.Range("J4").Value = 0.1601512700977636
$ws.Range("K4").Value = 1.137335566993244
$ws.Range("M4").Value = 0.384738083723235
$ws.Range("O4").Value = 3.876545931048128
$ws.Range("C5").Value = 0.03845195580052518
$ws.Range("D5").Value = 0.1492589943769715
$ws.Range("E5").Value = 0.1366304877349389
$ws.Range("F5").Value = 1.550806723961344
$ws.Range("G5").Value = 0.002479186706594536
$ws.Range("J5").Value = 0.160574259501221
$ws.Range("K5").Value = 1.100225151047198
$ws.Range("M5").Value = 0.3772673028289972
$ws.Range("O5").Value = 3.887934516381975
$ws.Range("C6").Value = 0.0382299448292116
$ws.Range("D6").Value = 0.1491686924377262
$ws.Range("E6").Value = 0.1366405314544963
$ws.Range("F6").Value = 1.55138941056552
$ws.Range("G6").Value = 0.002479335514321068
$ws.Range("J6").Value = 0.1606458455574682
$ws.Range("K6").Value = 1.094061600475442
$ws.Range("M6").Value = 0.3760279178164083
$ws.Range("O6").Value = 3.889868867689344
$ws.Range("C7").Value = 0.03977225839321363
$ws.Range("D7").Value = 0.1497991586832228
$ws.Range("E7").Value = 0.1365762118435345
$ws.Range("F7").Value = 1.547438336164717
$ws.Range("G7").Value = 0.002478312117620996
$ws.Range("J7").Value = 0.1601568842378818
$ws.Range("K7").Value = 1.136835177729751
$ws.Range("M7").Value = 0.3846372545042129
$ws.Range("O7").Value = 3.87669661889089
$ws.Range("C8").Value = 0.04661086479130461
$ws.Range("D8").Value = 0.1526717132156463
$ws.Range("E8").Value = 0.1364257712195496
$ws.Range("F8").Value = 1.532316830454064
$ws.Range("G8").Value = 0.002474030403928094
$ws.Range("J8").Value = 0.1581953810133889
$ws.Range("K8").Value = 1.325368903839205
$ws.Range("M8").Value = 0.422790167545152
$ws.Range("O8").Value = 3.824871197272984
$ws.Range("C9").Value = 0.06014127571451411
$ws.Range("D9").Value = 0.1586262487876979
$ws.Range("E9").Value = 0.1366093033082194
$ws.Range("F9").Value = 1.511002984422475
$ws.Range("G9").Value = 0.002466470628376524
$ws.Range("J9").Value = 0.1550548931259712
$ws.Range("K9").Value = 1.694191518976652
$ws.Range("M9").Value = 0.4981671073830825
$ws.Range("O9").Value = 3.745976384196638
$ws.Range("C10").Value = 0.07015866684930927
$ws.Range("D10").Value = 0.1631945441702811
$ws.Range("E10").Value = 0.137035366973798
$ws.Range("F10").Value = 1.500451212727739
$ws.Range("G10").Value = 0.002461421270909427
$ws.Range("J10").Value = 0.1531788942533083
$ws.Range("K10").Value = 1.964557874944603
$ws.Range("M10").Value = 0.5538752827900737
$ws.Range("O10").Value = 3.701956544225879
$ws.Range("C11").Value = 0.07473327775781513
$ws.Range("D11").Value = 0.1653143166412434
$ws.Range("E11").Value = 0.1372923388457181
$ws.Range("F11").Value = 1.496765169692225
$ws.Range("G11").Value = 0.002459232686670313
$ws.Range("J11").Value = 0.152419313293791
$ws.Range("K11").Value = 2.087411152086247
$ws.Range("M11").Value = 0.5792876321584117
$ws.Range("O11").Value = 3.684976967902969
$ws.Range("C12").Value = 0.07646814425078219
$ws.Range("D12").Value = 0.1661229546339058
$ws.Range("E12").Value = 0.1373987175955556
$ws.Range("F12").Value = 1.495529957893027
$ws.Range("G12").Value = 0.002458419429377214
$ws.Range("J12").Value = 0.1521451881140621
$ws.Range("K12").Value = 2.13391115510268
$ws.Range("M12").Value = 0.5889204528242402
$ws.Range("O12").Value = 3.67898663294099
$ws.Range("C13").Value = 0.07609439546843078
$ws.Range("D13").Value = 0.1659485373624108
$ws.Range("E13").Value = 0.1373754038734027
$ws.Range("F13").Value = 1.495788832886092
$ws.Range("G13").Value = 0.002458593890249892
$ws.Range("J13").Value = 0.1522036246169236
$ws.Range("K13").Value = 2.123897552494725
$ws.Range("M13").Value = 0.5868454269281642
$ws.Range("O13").Value = 3.680257188786584
$ws.Range("C14").Value = 0.07487595482037079
$ws.Range("D14").Value = 0.1653807253502464
$ws.Range("E14").Value = 0.1373009090016772
$ws.Range("F14").Value = 1.496660326755631
$ws.Range("G14").Value = 0.00245916546916635
$ws.Range("J14").Value = 0.1523964899395196
$ws.Range("K14").Value = 2.09123718526331
$ws.Range("M14").Value = 0.5800799377341406
$ws.Range("O14").Value = 3.684475323055636
$ws.Range("C15").Value = 0.0741299588823523
$ws.Range("D15").Value = 0.1650336940750918
$ws.Range("E15").Value = 0.1372564594405148
$ws.Range("F15").Value = 1.49721507086943
$ws.Range("G15").Value = 0.002459517595673436
$ws.Range("J15").Value = 0.1525163857401708
$ws.Range("K15").Value = 2.071228858922325
$ws.Range("M15").Value = 0.5759371333909229
$ws.Range("O15").Value = 3.687116329197011
$ws.Range("C16").Value = 0.06986006235752029
$ws.Range("D16").Value = 0.163056844613422
$ws.Range("E16").Value = 0.1370198430109539
$ws.Range("F16").Value = 1.500714556902651
$ws.Range("G16").Value = 0.002461566474749905
$ws.Range("J16").Value = 0.153230424302194
$ws.Range("K16").Value = 1.956526175969657
$ws.Range("M16").Value = 0.5522159077898152
$ws.Range("O16").Value = 3.703127611107845
$ws.Range("C17").Value = 0.06724515572939538
$ws.Range("D17").Value = 0.1618547312415473
$ws.Range("E17").Value = 0.1368908528014785
$ws.Range("F17").Value = 1.503147001556599
$ws.Range("G17").Value = 0.002462851104640824
$ws.Range("J17").Value = 0.1536925079603755
$ws.Range("K17").Value = 1.886123042632164
$ws.Range("M17").Value = 0.5376814465795974
$ws.Range("O17").Value = 3.713731087932558
$ws.Range("C18").Value = 0.06574279049955578
$ws.Range("D18").Value = 0.1611672296053968
$ws.Range("E18").Value = 0.1368226060675966
$ws.Range("F18").Value = 1.504650910293307
$ws.Range("G18").Value = 0.002463600196846511
$ws.Range("J18").Value = 0.1539671158493441
$ws.Range("K18").Value = 1.845616179567799
$ws.Range("M18").Value = 0.529328261491969
$ws.Range("O18").Value = 3.72011646952123
$ws.Range("C19").Value = 0.06523439984303536
$ws.Range("D19").Value = 0.1609351288534526
$ws.Range("E19").Value = 0.1368005203361733
$ws.Range("F19").Value = 1.505178099987447
$ws.Range("G19").Value = 0.002463855581799068
$ws.Range("J19").Value = 0.1540616092330929
$ws.Range("K19").Value = 1.831899109230847
$ws.Range("M19").Value = 0.5265011717461903
$ws.Range("O19").Value = 3.722327616504884
$ws.Range("C20").Value = 0.06752334504720636
$ws.Range("D20").Value = 0.1619822927727768
$ws.Range("E20").Value = 0.1369039688315112
$ws.Range("F20").Value = 1.502877211022238
$ws.Range("G20").Value = 0.002462713297971934
$ws.Range("J20").Value = 0.1536424043893021
$ws.Range("K20").Value = 1.893618920575193
$ws.Range("M20").Value = 0.5392279792658883
$ws.Range("O20").Value = 3.712572662394621
$ws.Range("C21").Value = 0.07523377069884418
$ws.Range("D21").Value = 0.1655473451468339
$ws.Range("E21").Value = 0.1373225439164933
$ws.Range("F21").Value = 1.496399985557403
$ws.Range("G21").Value = 0.002458997162076963
$ws.Range("J21").Value = 0.1523394738540524
$ws.Range("K21").Value = 2.100830936555553
$ws.Range("M21").Value = 0.5820668654813801
$ws.Range("O21").Value = 3.683224415498017
$ws.Range("C22").Value = 0.08028791533433832
$ws.Range("D22").Value = 0.1679118275230422
$ws.Range("E22").Value = 0.1376489628405935
$ws.Range("F22").Value = 1.493103112490871
$ws.Range("G22").Value = 0.002456658833862421
$ws.Range("J22").Value = 0.1515666991526281
$ws.Range("K22").Value = 2.236127575435262
$ws.Range("M22").Value = 0.6101209964134569
$ws.Range("O22").Value = 3.666605882761843
$ws.Range("C23").Value = 0.07758904913127651
$ws.Range("D23").Value = 0.1666467205227349
$ws.Range("E23").Value = 0.1374699146314136
$ws.Range("F23").Value = 1.494776899450656
$ws.Range("G23").Value = 0.002457898598331813
$ws.Range("J23").Value = 0.1519719297280311
$ws.Range("K23").Value = 2.163929634238286
$ws.Range("M23").Value = 0.5951429486643889
$ws.Range("O23").Value = 3.67524054703793
$ws.Range("C24").Value = 0.06739757260902479
$ws.Range("D24").Value = 0.1619246109583798
$ws.Range("E24").Value = 0.1368980206600661
$ws.Range("F24").Value = 1.502998854834146
$ws.Range("G24").Value = 0.00246277556748401
$ws.Range("J24").Value = 0.1536650283367642
$ws.Range("K24").Value = 1.890230131081466
$ws.Range("M24").Value = 0.5385287828686103
$ws.Range("O24").Value = 3.713095485518551
$ws.Range("C25").Value = 0.05646774150352485
$ws.Range("D25").Value = 0.1569812322759816
$ws.Range("E25").Value = 0.1365084747489895
$ws.Range("F25").Value = 1.515873716176316
$ws.Range("G25").Value = 0.002468426721602746
$ws.Range("J25").Value = 0.1558288005931311
$ws.Range("K25").Value = 1.594517412056291
$ws.Range("M25").Value = 0.4777170608361061
$ws.Range("O25").Value = 3.764876880213023
